$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.564.51"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.475.56"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.9554"
$ws.Range("E5").Value = "  +6.25%  "
$ws.Range("D6").Value = "'279.20"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.3650"
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").Value = "'40.02"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'1.060"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").Value = "'0.06673"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "'1.006"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'18.10"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'6.224"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'0.9551"
$ws.Range("E16").Value = "  +5.75%  "
$ws.Range("D17").Value = "'0.00001035"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "1.475.60"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'0.05956"
$ws.Range("D20").Value = "'69.94"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "'5.501"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").Value = "'14.47"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'11.08"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").Value = "'2.263"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "20.625.68"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'143.16"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("E27").Value = "  -6.39%  "
$ws.Range("D28").Value = "'17.28"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "1.636.40"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'114.06"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'3.962"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").Value = "'5.022"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "'0.8116"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "'0.07968"
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("D35").Value = "'1.514"
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("D36").Value = "'1.229"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("D37").Value = "'0.05845"
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("D38").Value = "'4.739"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'10.40"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "'0.9557"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("D42").Value = "'0.1881"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "'7.449"
$ws.Range("E43").Value = "  +6.94%  "
$ws.Range("D44").Value = "'0.5308"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.34"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.538"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'118.01"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").Value = "'0.5199"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'1.821"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'0.06483"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'0.9833"
$ws.Range("E51").Value = "  -1.10%  "
